$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 93
$ws1.Range("F3").Value = 65
$ws1.Range("F4").Value = 263
$ws1.Range("F5").Value = 151
$ws1.Range("F6").Value = 255
$ws1.Range("F7").Value = 206
$ws1.Range("F8").Value = 1940
$ws1.Range("F9").Value = 346
$ws1.Range("F10").Value = 4545
$ws1.Range("F11").Value = 71
$ws1.Range("F12").Value = 319

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2
$ws2.Range("F3").Value = 49
$ws2.Range("F4").Value = 6
$ws2.Range("F5").Value = 12

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2
$ws4.Range("F3").Value = 93
$ws4.Range("F4").Value = 65
$ws4.Range("F5").Value = 49
$ws4.Range("F6").Value = 263
$ws4.Range("F7").Value = 151
$ws4.Range("F8").Value = 255
$ws4.Range("F9").Value = 206
$ws4.Range("F10").Value = 6
$ws4.Range("F11").Value = 12
$ws4.Range("F12").Value = 1940
$ws4.Range("F13").Value = 346
$ws4.Range("F14").Value = 4545
$ws4.Range("F15").Value = 71
$ws4.Range("F16").Value = 319
